# Applies the "particion de atributos" (attribute/field partition) edit:
# adds "Has field"/"field" nodes and edges for PetController.pets and
# PetController.owners to the graph workbook ("nodes" and "Edges" sheets).

$wb = $excel.ActiveWorkbook

$wsNodes = $wb.Worksheets.Item("nodes")
$wsEdges = $wb.Worksheets.Item("Edges")

# --- New rows in "nodes" sheet (columns A..G), appended after row 38 ---
$nodesNew = @(
    @("pet2.org.springframework.samples.petclinic.owner.PetController.pets",  "org.springframework.samples.petclinic.owner", "PetController.pets",   "label", "field", "-", "pet2"),
    @("pet1.org.springframework.samples.petclinic.owner.PetController.pets",  "org.springframework.samples.petclinic.owner", "PetController.pets",   "label", "field", "-", "pet1"),
    @("pet1.org.springframework.samples.petclinic.owner.PetController.owners","org.springframework.samples.petclinic.owner", "PetController.owners", "label", "field", "-", "pet1")
)

$startRowNodes = 39
for ($i = 0; $i -lt $nodesNew.Length; $i++) {
    $r = $startRowNodes + $i
    $row = $nodesNew[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsNodes.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# --- New rows in "Edges" sheet (columns A..D), appended after row 23 ---
# Column D on existing rows carries a distinct (but visually identical)
# cell style (applyFont xf). Inserting a copy of the last existing row
# preserves that style on the new rows instead of recomputing a fresh one.
$edgesNew = @(
    @("pet2.org.springframework.samples.petclinic.owner.PetController", "pet2.org.springframework.samples.petclinic.owner.PetController.pets",   "Has field", "label"),
    @("pet1.org.springframework.samples.petclinic.owner.PetController", "pet1.org.springframework.samples.petclinic.owner.PetController.pets",   "Has field", "label"),
    @("pet1.org.springframework.samples.petclinic.owner.PetController", "pet1.org.springframework.samples.petclinic.owner.PetController.owners", "Has Field", "label")
)

$lastExistingRow = 23
$startRowEdges = 24
for ($i = 0; $i -lt $edgesNew.Length; $i++) {
    $r = $startRowEdges + $i
    $wsEdges.Rows.Item($lastExistingRow).Copy()
    $wsEdges.Rows.Item($r).Insert()

    $row = $edgesNew[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsEdges.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# --- Sheet view / selection changes ---
$wsNodes.Activate()
$wsNodes.Range("A29").Select()
$wsNodes.Application.ActiveWindow.Zoom = 115
$wsNodes.Range("G40").Select()

$wsEdges.Activate()
$wsEdges.Range("B26").Select()
